$p = $ppt.ActivePresentation

# Add a new slide 14 ("BEAM concurrency model") at the end of the deck.
# Duplicate an existing "Title and Content" slide (same placeholder
# naming / layout we need - title has a single run so re-stamping its
# text keeps a clean single <a:r>) and move the copy to the end, then
# overwrite its title & body text.
$base = $p.Slides.Item(10)
$ns = $base.Duplicate().Item(1)
$ns.MoveTo($p.Slides.Count)

# Title placeholder - clear first so no stale runs/formatting bleed
# through from the duplicated slide's title text.
$title = $ns.Shapes.Item(1).TextFrame.TextRange
$title.Text = ""
$title.Text = "BEAM concurrency model"

# Body placeholder - same trick: clear to empty before writing the new
# paragraphs so each one starts as a single, cleanly-formatted run
# instead of inheriting the old paragraph's pPr/runs positionally.
$body = $ns.Shapes.Item(2).TextFrame.TextRange
$body.Text = ""
$body.Text = "BEAM VM has its own scheduler`rRuns lightweight processes instead of OS threads`rNo shared state between processes`rCommunication is done via message passing`rIsolation means:`rBetter performance (no locks)`rFault tolerance`rScalability – up to ~268M running processes"

$body.Paragraphs(2,1).IndentLevel = 2
$body.Paragraphs(3,1).IndentLevel = 2
$body.Paragraphs(4,1).IndentLevel = 2
$body.Paragraphs(6,1).IndentLevel = 2
$body.Paragraphs(7,1).IndentLevel = 2
$body.Paragraphs(8,1).IndentLevel = 2
